# Update Wnt5a-Ror1 LR-pairs sheet with refreshed TPM-derived NATMI values.
# The sending/target cluster set changed from {FAPs -> ECs,FAPs,MuSCs,Neutrophils}
# to {ECs,FAPs} x {ECs,FAPs,MuSCs}: the "Neutrophils" target row is gone, two new
# "ECs"-sending rows were added, and all numeric NATMI metrics were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs","Wnt5a","Ror1","ECs",2,1,0.125615,0.25123,0.02647478672532295,0.01780700335556722,2,1,1.592395,3.18479,0.1081406801663129,0.08970879337942068,0.200028697925,0.8001147916999999,0.002863001443734495,0.00159744478473123),
    @("ECs","Wnt5a","Ror1","FAPs",2,1,0.125615,0.25123,0.02647478672532295,0.01780700335556722,3,1,6.050990666666666,18.152972,0.4109270918166311,0.5113307986932919,0.7600951925933334,4.56057115556,0.01087920711550251,0.009105269248136316),
    @("ECs","Wnt5a","Ror1","MuSCs",2,1,0.125615,0.25123,0.02647478672532295,0.01780700335556722,2,1,7.0818315,14.163663,0.4809322280170559,0.3989604079272874,0.8895842638725,3.55833705549,0.01273257816608594,0.007104289322699674),
    @("FAPs","Wnt5a","Ror1","ECs",3,1,4.619088000000001,13.857264,0.9735252132746771,0.9821929966444328,2,1,1.592395,3.18479,0.1081406801663129,0.08970879337942068,7.35541263576,44.13247581456,0.1052776787225784,0.08811134859468946),
    @("FAPs","Wnt5a","Ror1","FAPs",3,1,4.619088000000001,13.857264,0.9735252132746771,0.9821929966444328,3,1,6.050990666666666,18.152972,0.4109270918166311,0.5113307986932919,27.950058376512,251.550525388608,0.4000478847011286,0.5022255294451556),
    @("FAPs","Wnt5a","Ror1","MuSCs",3,1,4.619088000000001,13.857264,0.9735252132746771,0.9821929966444328,2,1,7.0818315,14.163663,0.4809322280170559,0.3989604079272874,32.711602899672,196.269617398032,0.46819964985097,0.3918561186045878)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $addr = $cols[$c] + $excelRow
        $ws.Range($addr).Value = $rowVals[$c]
    }
}
